# Add new column 'Servised by' to Card24
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card24")

# New header cell O1: copy the formatting of the neighbouring header cell
# (N1, "Correction") straight into O1 - Copy(Destination) carries the
# header style (bold / centered / bordered) across verbatim - then set
# its own text.
$ws.Range("N1").Copy($ws.Range("O1")) | Out-Null
$ws.Range("O1").Value = "Servised by"

# L2 should actually hold "nan" like the rest of row 2, instead of being
# blank.
$ws.Range("L2").Value = "nan"

# Materialize O2:O12 as existing (but empty) cells now that column O is
# part of the sheet's used range, without altering their (default)
# formatting (explicitly (re)asserting "no border" touches the cell
# without creating/leaving behind any new style).
for ($r = 2; $r -le 12; $r++) {
    $ws.Cells.Item($r, 15).Borders.LineStyle = -4142
}
